# Perbaikan Antrian Device Presensi
# Update attendance times and statuses, then refresh the summary counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose "Hadir" time shifts from 18:40 -> 19:47
$ws.Range("E5").Value = "19:47"
$ws.Range("E6").Value = "19:47"
$ws.Range("E22").Value = "19:47"
$ws.Range("E26").Value = "19:47"
$ws.Range("E28").Value = "19:47"
$ws.Range("E35").Value = "19:47"

# Rows whose status flips from "Hadir" (18:43) to "Alpha" (00:00)
$ws.Range("D8").Value = "Alpha"
$ws.Range("E8").Value = "00:00"

$ws.Range("D19").Value = "Alpha"
$ws.Range("E19").Value = "00:00"

$ws.Range("D33").Value = "Alpha"
$ws.Range("E33").Value = "00:00"

# Update the summary counts at the bottom of the report
$ws.Range("A38").Value = "Hadir: 6"
$ws.Range("A41").Value = "Alpha: 28"
